$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.379.28'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.846.65'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('D4').Value = '''0.9993'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''241.29'
$ws.Range('D6').Value = '''0.6333'
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''0.07562'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').Value = '''0.2958'
$ws.Range('E9').Value = '  -1.50%  '
$ws.Range('D10').Value = '''24.77'
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('D11').Value = '''0.07743'
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').Value = '''4.994'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').Value = '''0.6822'
$ws.Range('E13').Value = '  -1.09%  '
$ws.Range('D14').Value = '''83.02'
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').Value = '''0.000009962'
$ws.Range('E15').Value = '  +2.74%  '
$ws.Range('D16').Value = '''6.140'
$ws.Range('E16').Value = '  -2.70%  '
$ws.Range('D17').Value = '29.403.11'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').Value = '''230.30'
$ws.Range('E18').Value = '  -3.78%  '
$ws.Range('E19').Value = '  -1.28%  '
$ws.Range('D20').Value = '''1.000'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '''7.551'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('E23').Value = '  +232.74%  '
$ws.Range('D24').Value = '''16.60'
$ws.Range('E24').Value = '  +170.27%  '
$ws.Range('D25').Value = '''156.47'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').Value = '''0.1399'
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('D27').Value = '''8.388'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('D28').Value = '''17.67'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').Value = '''2.733'
$ws.Range('E29').Value = '  +173.48%  '
$ws.Range('D30').Value = '''1.473'
$ws.Range('E30').Value = '  -1.04%  '
$ws.Range('D31').Value = '''0.05714'
$ws.Range('E31').Value = '  -3.04%  '
$ws.Range('D32').Value = '''1.254'
$ws.Range('E32').Value = '  -2.46%  '
$ws.Range('D33').Value = '''4.126'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('D34').Value = '''4.009'
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('D35').Value = '''1.844'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('D36').Value = '''1.155'
$ws.Range('E36').Value = '  -2.59%  '
$ws.Range('D37').Value = '''0.7165'
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('D38').Value = '''2.595'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = '1.244.04'
$ws.Range('E39').Value = '  +2.07%  '
$ws.Range('D40').Value = '''2.801'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').Value = '''0.01808'
$ws.Range('E41').Value = '  +1.55%  '
$ws.Range('E42').Value = '  +265.21%  '
$ws.Range('D43').Value = '''0.9014'
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('D44').Value = '''1.001'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = '''101.87'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('D47').Value = '''7.048'
$ws.Range('E47').Value = '  -5.92%  '
$ws.Range('D48').Value = '''9.108'
$ws.Range('E48').Value = '  -1.05%  '
$ws.Range('D49').Value = '''0.4017'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').Value = '''1.701'
$ws.Range('D51').Value = '''0.1123'
$ws.Range('E51').Value = '  -0.89%  '
